$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) New paragraph: "LSTM understanding the Number of Parameters" + ":"
#    (kept as two separate runs, matching the source edit)
# ---------------------------------------------------------------------
$end = $d.Content.End
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.InsertAfter("LSTM understanding the Number of Parameters")

$afterText1 = $d.Content.End
$splitRng = $d.Range($afterText1, $afterText1)
$splitRng.InsertParagraphAfter()
$afterSplit1 = $d.Content.End
$colonRng = $d.Range($afterSplit1 - 1, $afterSplit1 - 1)
$colonRng.InsertAfter(":")
# Delete the temporary paragraph mark that separated the two runs so the
# colon ends up as its own run inside the same paragraph as the title text.
$joinRng = $d.Range($afterText1 - 1, $afterText1)
$joinRng.Delete()

# ---------------------------------------------------------------------
# 2) New paragraph containing the hyperlink to the Colab notebook
# ---------------------------------------------------------------------
$end2 = $d.Content.End
$rng2 = $d.Range($end2, $end2)
$rng2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
# Insert a placeholder character so the hyperlink is not anchored at an
# entirely-empty paragraph (which the engine mis-serializes); it is
# removed again immediately after the hyperlink is created.
$p2.Range.InsertAfter("X")
$hStart = $p2.Range.Start
$hRng = $d.Range($hStart, $hStart)
$d.Hyperlinks.Add($hRng, "https://colab.research.google.com/github/kmkarakaya/ML_tutorials/blob/master/LSTM_Understanding_the_Number_of_Parameters.ipynb") | Out-Null
$p2End = $p2.Range.End
$placeholder2 = $d.Range($p2End - 2, $p2End - 1)
$placeholder2.Delete()

# ---------------------------------------------------------------------
# 3) Trailing empty paragraph holding the "_GoBack" bookmark Word leaves
#    at the last edited location.
# ---------------------------------------------------------------------
$end3 = $d.Content.End
$rng3 = $d.Range($end3, $end3)
$rng3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.InsertAfter("X")
$bPos = $p3.Range.Start
$bRng = $d.Range($bPos, $bPos)
$d.Bookmarks.Add("_GoBack", $bRng)
$placeholder3 = $d.Range($bPos, $bPos + 1)
$placeholder3.Delete()

# ---------------------------------------------------------------------
# 4) Register the (now used) built-in "Hyperlink" character style with
#    the same definition Word writes out when AutoFormat creates one.
# ---------------------------------------------------------------------
$hlStyle = $d.Styles("Hyperlink")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.QuickStyle = $false

Write-Host "Paragraphs: " $d.Paragraphs.Count
